$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "https://github.com/Kabesa/images/blob/master/sportsCategories/vetements.png?raw=true"
$ws.Range("C3").Value = "https://github.com/Kabesa/images/blob/master/sportsCategories/bicycle.png?raw=true"
$ws.Range("C4").Value = "https://github.com/Kabesa/images/blob/master/sportsCategories/ski.png?raw=true"
$ws.Range("C5").Value = "https://github.com/Kabesa/images/blob/master/sportsCategories/golf.png?raw=true"
$ws.Range("C6").Value = "https://github.com/Kabesa/images/blob/master/sportsCategories/camping.png?raw=true"

$ws.Range("C3").Select()
